$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the birth date column (F) for several rows with real dates
# (values are Excel date serial numbers, matching the existing cell format)
$ws.Range("F2").Value = 36526
$ws.Range("F3").Value = 36161
$ws.Range("F4").Value = 37622
$ws.Range("F6").Value = 37257
$ws.Range("F7").Value = 35431
$ws.Range("F8").Value = 35065
$ws.Range("F9").Value = 38718

# Move the active selection to F6 (was D10)
$ws.Range("F6").Select()
